# feat: add 2022-Q3 data
#
# - Duplicate the "2022-Q2" sheet, place the copy right before it, rename
#   the copy to "2022-Q3" (this naturally renumbers/shifts every sheet
#   after it, exactly like the target workbook layout).
# - Update the fund snapshot numbers on the new "2022-Q3" sheet.
# - Update the "总计" (summary) sheet: the quarters shift down one row and
#   a new row is added for "2021-Q1" (2020-Q4 moves to row 7).

$wb = $excel.ActiveWorkbook

# ---------------------------------------------------------------------
# 1. Create the new "2022-Q3" sheet as a copy of "2022-Q2", inserted
#    immediately before it.
# ---------------------------------------------------------------------
$q2 = $wb.Worksheets.Item("2022-Q2")
$q2.Copy($q2, [System.Reflection.Missing]::Value)
$q3 = $wb.Worksheets.Item("2022-Q2 (2)")
$q3.Name = "2022-Q3"

# Helper: write a value as TEXT (matching the source file, which stores
# these numbers as inline strings, not numeric cells) without leaving a
# stray "quote prefix" style behind.
function Set-TextValue($range, [string]$text) {
    $range.Value = "'" + $text
    $range.Style = "Normal"
}

Set-TextValue $q3.Range("D2") "13.65"
Set-TextValue $q3.Range("E2") "73.02"
Set-TextValue $q3.Range("F2") "2.95"
Set-TextValue $q3.Range("G2") "0.4027"

# ---------------------------------------------------------------------
# 2. Update the "总计" summary sheet.
# ---------------------------------------------------------------------
$summary = $wb.Worksheets.Item("总计")

$summary.Range("B2").Value = "2022-Q3"
$summary.Range("C2").Value = 1
$summary.Range("D2").Value = 0.4

$summary.Range("B3").Value = "2022-Q2"
$summary.Range("C3").Value = 1
$summary.Range("D3").Value = 0.25

$summary.Range("B4").Value = "2021-Q4"
$summary.Range("C4").Value = 2
$summary.Range("D4").Value = 0.66

$summary.Range("B5").Value = "2021-Q2"
$summary.Range("C5").Value = 1
$summary.Range("D5").Value = 0.01

# New row 6 ("2021-Q1") -- copy formatting from row 6's A-cell (which
# currently holds the soon-to-move-down "2020-Q4" row) before overwriting
# its value, so the style carries over correctly.
$summary.Range("B6").Value = "2021-Q1"
$summary.Range("C6").Value = 1
$summary.Range("D6").Value = 0.01

# New row 7 ("2020-Q4", moved down from row 6). Clone A6's style onto A7
# first (A7 doesn't exist yet so it has no formatting of its own), then
# set the real values for the row.
$summary.Range("A6").Copy($summary.Range("A7"))
$summary.Range("A7").Value = 5
$summary.Range("B7").Value = "2020-Q4"
$summary.Range("C7").Value = 3
$summary.Range("D7").Value = 0.01

Write-Host "2022-Q3 sheet added and 总计 summary updated"
